# Add two new columns I (I0) and J (IF) to the sheet, mirroring the
# existing H (IP) column layout: same header style, same plain numeric
# body cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1): copy the format of the existing H1 header cell so
# the new header cells pick up the same style (bold font, border,
# centered/top alignment) instead of creating brand new style entries.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Body values (rows 2-23) ---
# I column: "I0" is 1 for every game except row 18, where it is 4.
# J column: "IF" = H (IP) + I (I0) - 1 for every row.
$i0 = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 4;
    19 = 1; 20 = 1; 21 = 1; 22 = 1; 23 = 1
}

for ($row = 2; $row -le 23; $row++) {
    $ip = $ws.Cells.Item($row, 8).Value2
    $iVal = $i0[$row]
    $jVal = $ip + $iVal - 1

    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
